# Update "想去人数" (want-to-go headcount) figures in column F
# across all four sheets, per the 456a3b4 gh-pages data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 532
$ws.Range("F4").Value = 628
$ws.Range("F6").Value = 458
$ws.Range("F8").Value = 2119
$ws.Range("F9").Value = 868
$ws.Range("F10").Value = 827
$ws.Range("F12").Value = 74
$ws.Range("F14").Value = 320
$ws.Range("F15").Value = 99
$ws.Range("F18").Value = 30
$ws.Range("F19").Value = 1696
$ws.Range("F27").Value = 520
$ws.Range("F28").Value = 349
$ws.Range("F29").Value = 587
$ws.Range("F30").Value = 411
$ws.Range("F31").Value = 2350
$ws.Range("F33").Value = 87
$ws.Range("F35").Value = 596
$ws.Range("F36").Value = 469
$ws.Range("F37").Value = 181
$ws.Range("F38").Value = 915
$ws.Range("F41").Value = 427
$ws.Range("F42").Value = 390

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 51
$ws.Range("F14").Value = 74
$ws.Range("F21").Value = 87
$ws.Range("F22").Value = 117
$ws.Range("F23").Value = 98
$ws.Range("F24").Value = 433

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2917
$ws.Range("F6").Value = 306

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 532
$ws.Range("F7").Value = 628
$ws.Range("F10").Value = 458
$ws.Range("F11").Value = 868
$ws.Range("F12").Value = 827
$ws.Range("F14").Value = 74
$ws.Range("F15").Value = 320
$ws.Range("F16").Value = 99
$ws.Range("F21").Value = 30
$ws.Range("F22").Value = 306
$ws.Range("F23").Value = 1696
$ws.Range("F27").Value = 51
$ws.Range("F30").Value = 74
$ws.Range("F34").Value = 520
$ws.Range("F35").Value = 349
$ws.Range("F36").Value = 587
$ws.Range("F37").Value = 411
$ws.Range("F38").Value = 87
$ws.Range("F40").Value = 469
$ws.Range("F41").Value = 181
$ws.Range("F42").Value = 915
$ws.Range("F43").Value = 87
$ws.Range("F44").Value = 98
$ws.Range("F45").Value = 433
$ws.Range("F48").Value = 427
$ws.Range("F49").Value = 390

